$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 31: restrict E31 date format (s=3 -> s=2) and I31 type (boolean -> number) ---
$ws.Cells.Item(31, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31, 9).Value = 1

# --- Append new rows 32-59 ---
# Row 32
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 5).Value = 45428
$ws.Cells.Item(32, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(32, 6).Value = 22
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 1

# Row 33
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 5).Value = 45428
$ws.Cells.Item(33, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 1

# Row 34
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 5).Value = 45429
$ws.Cells.Item(34, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 1

# Row 35
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 5).Value = 45429
$ws.Cells.Item(35, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 1

# Row 36
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 5).Value = 45429
$ws.Cells.Item(36, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 1

# Row 37
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 5).Value = 45429
$ws.Cells.Item(37, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 1

# Row 38
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = 2
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = "aiaiai"
$ws.Cells.Item(38, 5).Value = 45429
$ws.Cells.Item(38, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 1

# Row 39
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = "ala"
$ws.Cells.Item(39, 5).Value = 45429
$ws.Cells.Item(39, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 1

# Row 40
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = 2
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 5).Value = 45429
$ws.Cells.Item(40, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 1

# Row 41
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 2
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = "aju"
$ws.Cells.Item(41, 5).Value = 45429
$ws.Cells.Item(41, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 1

# Row 42
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = "aju"
$ws.Cells.Item(42, 5).Value = 45429
$ws.Cells.Item(42, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 1

# Row 43
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 2
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = "aju"
$ws.Cells.Item(43, 5).Value = 45429
$ws.Cells.Item(43, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 1

# Row 44
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 2
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = "aju"
$ws.Cells.Item(44, 5).Value = 45429
$ws.Cells.Item(44, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 1

# Row 45
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 2
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = "aju"
$ws.Cells.Item(45, 5).Value = 45429
$ws.Cells.Item(45, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 1

# Row 46
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = "aju"
$ws.Cells.Item(46, 5).Value = 45429
$ws.Cells.Item(46, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 1

# Row 47
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = 2
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(47, 4).Value = "aju"
$ws.Cells.Item(47, 5).Value = 45429
$ws.Cells.Item(47, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 1

# Row 48
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 2
$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(48, 4).Value = "aju"
$ws.Cells.Item(48, 5).Value = 45429
$ws.Cells.Item(48, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 1

# Row 49
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = 2
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = "aju"
$ws.Cells.Item(49, 5).Value = 45429
$ws.Cells.Item(49, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 9).Value = 1

# Row 50
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = 2
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(50, 4).Value = "aju"
$ws.Cells.Item(50, 5).Value = 45429
$ws.Cells.Item(50, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 1

# Row 51
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = 2
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(51, 4).Value = "aju"
$ws.Cells.Item(51, 5).Value = 45429
$ws.Cells.Item(51, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 1

# Row 52
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 3).Value = 1
$ws.Cells.Item(52, 4).Value = "aju"
$ws.Cells.Item(52, 5).Value = 45429
$ws.Cells.Item(52, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 1

# Row 53
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = 2
$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(53, 4).Value = "aju"
$ws.Cells.Item(53, 5).Value = 45429
$ws.Cells.Item(53, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 1

# Row 54
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = 2
$ws.Cells.Item(54, 3).Value = 1
$ws.Cells.Item(54, 4).Value = "aju"
$ws.Cells.Item(54, 5).Value = 45429
$ws.Cells.Item(54, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 1

# Row 55
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = 2
$ws.Cells.Item(55, 3).Value = 1
$ws.Cells.Item(55, 4).Value = "aju"
$ws.Cells.Item(55, 5).Value = 45429
$ws.Cells.Item(55, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 1

# Row 56
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = 2
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(56, 4).Value = "aju"
$ws.Cells.Item(56, 5).Value = 45429
$ws.Cells.Item(56, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 1

# Row 57
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = 2
$ws.Cells.Item(57, 3).Value = 1
$ws.Cells.Item(57, 4).Value = "ajo"
$ws.Cells.Item(57, 5).Value = 45429
$ws.Cells.Item(57, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 1

# Row 58
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = 2
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(58, 4).Value = "ajo"
$ws.Cells.Item(58, 5).Value = 45429
$ws.Cells.Item(58, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 1

# Row 59
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = 2
$ws.Cells.Item(59, 3).Value = 1
$ws.Cells.Item(59, 4).Value = "ajo"
$ws.Cells.Item(59, 5).Value = 45429
$ws.Cells.Item(59, 5).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 9).Value = $true

